$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value (applied uniformly across columns J:AS)
$updates = @{
    100 = 3991.635443
    101 = 151205.7682
    102 = 406914.786
    103 = 31250.14823
    104 = 3070.488802
    105 = 73827.77658999999
    106 = 44760.62646
    107 = 84212.49442
    114 = 1715.427119
    115 = 385707.7218
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $range = $ws.Range("J$row`:AS$row")
    $range.Value = $value
}
